$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 63
$ws.Range("D2").Value = 126
$ws.Range("E2").Value = 33.33333333333333
$ws.Range("F2").Value = 66.66666666666666

$ws.Range("C3").Value = 147
$ws.Range("D3").Value = 87
$ws.Range("E3").Value = 62.82051282051282
$ws.Range("F3").Value = 37.17948717948718

$ws.Range("C4").Value = 116
$ws.Range("E4").Value = 46.21513944223107
$ws.Range("F4").Value = 53.78486055776892

$ws.Range("E7").Value = 47.57804090419806
$ws.Range("F7").Value = 52.42195909580194
